# Updated symbol list on Tue Feb  7 13:55:34 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns on the crypto tracker sheet.
# Values are stored as literal text (inlineStr/shared-string), matching the
# source sheet's layout, so we force the Text number format before writing
# and reset the style afterwards to avoid leaving stray cell formatting.

function Set-TextCellValue {
    param($Worksheet, $Address, $Text)
    $cell = $Worksheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCellValue $ws "D2" "329.24"
Set-TextCellValue $ws "E2" "0.53%"
Set-TextCellValue $ws "D3" "44.28"
Set-TextCellValue $ws "E3" "0.99%"
Set-TextCellValue $ws "D4" "5.517"
Set-TextCellValue $ws "E4" "-0.83%"
Set-TextCellValue $ws "D5" "0.08019"
Set-TextCellValue $ws "E5" "-0.32%"
Set-TextCellValue $ws "D6" "2.071"
Set-TextCellValue $ws "E6" "8.57%"
Set-TextCellValue $ws "D7" "2.627"
Set-TextCellValue $ws "E7" "3.41%"
Set-TextCellValue $ws "D8" "0.9555"
Set-TextCellValue $ws "E8" "1.02%"
Set-TextCellValue $ws "D9" "0.1144"
Set-TextCellValue $ws "E9" "-1.92%"
Set-TextCellValue $ws "D10" "0.1892"
Set-TextCellValue $ws "E10" "2.74%"
Set-TextCellValue $ws "D11" "10.18"
Set-TextCellValue $ws "E11" "6.35%"
Set-TextCellValue $ws "D12" "0.09869"
Set-TextCellValue $ws "E12" "1.22%"
Set-TextCellValue $ws "D13" "0.04868"
Set-TextCellValue $ws "E13" "10.72%"
Set-TextCellValue $ws "E14" "-0.57%"
Set-TextCellValue $ws "D15" "0.001280"
Set-TextCellValue $ws "E15" "-0.22%"
Set-TextCellValue $ws "D16" "0.04089"
Set-TextCellValue $ws "E16" "-2.92%"
Set-TextCellValue $ws "D17" "0.006142"
Set-TextCellValue $ws "E17" "3.11%"
Set-TextCellValue $ws "E18" "-1.06%"
Set-TextCellValue $ws "D19" "4.398"
Set-TextCellValue $ws "E19" "2.55%"
Set-TextCellValue $ws "D20" "0.3405"
Set-TextCellValue $ws "E20" "-2.56%"
Set-TextCellValue $ws "D21" "0.1382"
Set-TextCellValue $ws "E21" "0.29%"
Set-TextCellValue $ws "D22" "0.2578"
Set-TextCellValue $ws "E22" "2.75%"
Set-TextCellValue $ws "D23" "0.001302"
Set-TextCellValue $ws "E23" "4.37%"
Set-TextCellValue $ws "E24" "0.57%"
Set-TextCellValue $ws "E25" "-6.47%"
Set-TextCellValue $ws "D26" "0.0003745"
Set-TextCellValue $ws "E26" "-6.30%"
Set-TextCellValue $ws "D38" "0.02583"
Set-TextCellValue $ws "E38" "-2.20%"
Set-TextCellValue $ws "D39" "0.05775"
Set-TextCellValue $ws "E39" "5.16%"
Set-TextCellValue $ws "D40" "0.007578"
Set-TextCellValue $ws "E40" "-0.23%"
Set-TextCellValue $ws "D41" "0.1402"
Set-TextCellValue $ws "E41" "0.41%"
Set-TextCellValue $ws "D42" "0.007338"
Set-TextCellValue $ws "E42" "-8.87%"
Set-TextCellValue $ws "D43" "0.002008"
Set-TextCellValue $ws "E43" "-0.19%"
Set-TextCellValue $ws "D44" "0.009048"
Set-TextCellValue $ws "E44" "4.96%"
Set-TextCellValue $ws "D45" "0.00007012"
Set-TextCellValue $ws "E45" "1.33%"
Set-TextCellValue $ws "E46" "-0.17%"
Set-TextCellValue $ws "D48" "0.003499"
Set-TextCellValue $ws "E48" "53.84%"
Set-TextCellValue $ws "E49" "-33.25%"
Set-TextCellValue $ws "D50" "0.00002101"
Set-TextCellValue $ws "E50" "-0.17%"
Set-TextCellValue $ws "D51" "0.0002001"
Set-TextCellValue $ws "E51" "-0.17%"
